$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.346648216247559
$ws.Range("B1").Value = 2.63078498840332
$ws.Range("C1").Value = 2.261423110961914
$ws.Range("D1").Value = 2.398743867874146
$ws.Range("E1").Value = 2.861247539520264
